$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4
$ws.Range("D4").Value = 44300
$ws.Range("M4").Value = 100
$ws.Range("N4").Value = 15000
$ws.Range("O4").Value = 15000
$ws.Range("P4").Value = 15000
$ws.Range("S4").Value = 2143

# Row 5
$ws.Range("D5").Value = 44300
$ws.Range("N5").Value = 12000
$ws.Range("O5").Value = 12000
$ws.Range("P5").Value = 12000
$ws.Range("S5").Value = 1714

# Row 6
$ws.Range("D6").Value = 44302
$ws.Range("M6").Value = 50
$ws.Range("R6").Value = "Región Metropolitana"

# Row 7
$ws.Range("D7").Value = 44302
$ws.Range("M7").Value = 30
$ws.Range("R7").Value = "Región Metropolitana"

# Row 8
$ws.Range("D8").Value = 44301
$ws.Range("M8").Value = 100
$ws.Range("N8").Value = 14000
$ws.Range("O8").Value = 14000
$ws.Range("P8").Value = 14000
$ws.Range("S8").Value = 2000

# Row 9
$ws.Range("D9").Value = 44301
$ws.Range("M9").Value = 80
$ws.Range("N9").Value = 12000
$ws.Range("O9").Value = 12000
$ws.Range("P9").Value = 12000
$ws.Range("S9").Value = 1714

# Row 10
$ws.Range("D10").Value = 44322
$ws.Range("M10").Value = 45
$ws.Range("N10").Value = 12000
$ws.Range("O10").Value = 12000
$ws.Range("P10").Value = 12000
$ws.Range("S10").Value = 1714

# Row 11
$ws.Range("D11").Value = 44322
$ws.Range("N11").Value = 8000
$ws.Range("O11").Value = 8000
$ws.Range("P11").Value = 8000
$ws.Range("S11").Value = 1143

# Row 12
$ws.Range("D12").Value = 44292
$ws.Range("M12").Value = 25
$ws.Range("N12").Value = 16000
$ws.Range("O12").Value = 16000
$ws.Range("P12").Value = 16000
$ws.Range("S12").Value = 2286

# Row 13
$ws.Range("D13").Value = 44292
$ws.Range("N13").Value = 15000
$ws.Range("O13").Value = 15000
$ws.Range("P13").Value = 15000
$ws.Range("S13").Value = 2143

# Row 14
$ws.Range("D14").Value = 44299
$ws.Range("M14").Value = 80
$ws.Range("N14").Value = 15000
$ws.Range("O14").Value = 15000
$ws.Range("P14").Value = 15000
$ws.Range("R14").Value = "Provincia de Santiago"
$ws.Range("S14").Value = 2143

# Row 15
$ws.Range("D15").Value = 44299
$ws.Range("M15").Value = 75
$ws.Range("R15").Value = "Provincia de Santiago"
